# Rename header columns in the shared-string table:
#   "<Name>_old" -> "<Name>_FV2410"
#   "<Name>_new" -> "<Name>_FV2504"
# then wrap the used range in an Excel Table ("Table1") using those same
# (renamed) header values, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J carry the "_old" -> "_FV2410" headers, columns L..U carry the
# "_new" -> "_FV2504" headers; column K is the untouched "diff" column.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $oldCol = 1 + $i        # A..J
    $newCol = 12 + $i       # L..U
    $ws.Cells.Item(1, $oldCol).Value = "$($baseNames[$i])_FV2410"
    $ws.Cells.Item(1, $newCol).Value = "$($baseNames[$i])_FV2504"
}

# Wrap the data region in a proper Excel Table named "Table1".
$range = $ws.Range("A1:U65")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# Freeze the header row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
